$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.827.76"
$ws.Range("E2").Value = "  +3.09%  "

$ws.Range("D3").Value = "3.439.46"
$ws.Range("E3").Value = "  +1.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.59"
$ws.Range("E5").Value = "  +2.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.46"
$ws.Range("E6").Value = "  +4.58%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.476"
$ws.Range("E8").Value = "  +0.83%  "

$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  +2.14%  "

$ws.Range("E11").Value = "  +2.22%  "

$ws.Range("D12").Value = "4.033.07"
$ws.Range("E12").Value = "  +1.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.37"
$ws.Range("E13").Value = "  +5.81%  "

$ws.Range("D15").Value = "3.440.02"
$ws.Range("E15").Value = "  +0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("D17").Value = "62.853.71"
$ws.Range("E17").Value = "  +2.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.21"
$ws.Range("E18").Value = "  +2.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.22"
$ws.Range("E19").Value = "  +5.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.27"
$ws.Range("E20").Value = "  +4.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.25"
$ws.Range("E21").Value = "  +3.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.38"
$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.560"
$ws.Range("E23").Value = "  +2.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000118"
$ws.Range("E25").Value = "  +4.62%  "

$ws.Range("D26").Value = "3.584.29"
$ws.Range("E26").Value = "  +1.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.191"
$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.66"
$ws.Range("E28").Value = "  +6.89%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.17"
$ws.Range("E30").Value = "  +2.99%  "

$ws.Range("E31").Value = "  +1.31%  "

$ws.Range("E32").Value = "  +5.92%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.69"
$ws.Range("E34").Value = "  +2.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  +7.02%  "

$ws.Range("E36").Value = "  +2.57%  "

$ws.Range("E37").Value = "  +8.45%  "

$ws.Range("E38").Value = "  +1.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "30.49"
$ws.Range("E39").Value = "  +18.94%  "

$ws.Range("D40").Value = "3.475.51"
$ws.Range("E40").Value = "  +1.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0765"
$ws.Range("E41").Value = "  +0.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.791"
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.88"
$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("E44").Value = "  +3.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("E45").Value = "  +5.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.19"
$ws.Range("E46").Value = "  +7.88%  "

$ws.Range("D47").Value = "2.524.04"
$ws.Range("E47").Value = "  +3.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.50"
$ws.Range("E48").Value = "  +3.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.72"
$ws.Range("E49").Value = "  +1.94%  "

$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.18"
$ws.Range("E51").Value = "  +4.51%  "
